$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text used for the newly inserted "Possible_Problem" rows (column C), identical to the
# text already used in rows 7/10 (shared string index 14 in the original file).
$possibleProblemLabel = "Possible_Problem"
$wrapText = "Possible_Problem:35% Power Window Regulator`n25% Power Window Motor`n20% Window Switch`n15% Door Jam Wiring`n5% Restricted window tracks"

# --- Insert new row before row 17 ---------------------------------------------------
# Node column (A) repeats the node label already present in the row right above it.
$nodeA = $ws.Range("A16").Value2
$ws.Rows("17:17").Insert()
$ws.Range("A17").Value = $nodeA
$ws.Range("B17").Value = $possibleProblemLabel
$ws.Range("C17").Value = $wrapText
$ws.Range("C17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 288

# --- Insert new row before row 20 ---------------------------------------------------
$nodeA = $ws.Range("A19").Value2
$ws.Rows("20:20").Insert()
$ws.Range("A20").Value = $nodeA
$ws.Range("B20").Value = $possibleProblemLabel
$ws.Range("C20").Value = $wrapText
$ws.Range("C20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 288

# --- Insert new row before row 22 ---------------------------------------------------
$nodeA = $ws.Range("A23").Value2
$ws.Rows("22:22").Insert()
$ws.Range("A22").Value = $nodeA
$ws.Range("B22").Value = $possibleProblemLabel
$ws.Range("C22").Value = $wrapText
$ws.Range("C22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 288

# Restore default (A1) selection instead of the stale B3 selection left in the source file.
$ws.Range("A1").Select() | Out-Null
